$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D in the original sheet stores prices as plain *text* (inline strings),
# e.g. "26.193.96" or "209.55". Assigning a plain numeric-looking string to .Value
# would make Excel auto-convert it to a real number, losing the text formatting/value
# fidelity (e.g. "209.56" turning into 209.56 as a float, or trailing zeros being lost).
# To preserve the original text semantics we force the cell format to Text ("@") right
# before writing values that look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.183.94'
$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.580.64'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.56'
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("E6").Value = '  -3.28%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0609'
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.803.64'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.614.15'
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.05'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.35'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.195.02'
$ws.Range("E17").Value = '  -1.90%  '
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.28'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '207.11'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.25'
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.45'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.99'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0505'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.14'
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.276.40'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.46'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.610'
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.817'
$ws.Range("E39").Value = '  -1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  -11.61%  '
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.13'
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.765'
$ws.Range("E43").Value = '  -2.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.19'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.716.50'
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.02'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.46'
$ws.Range("E51").Value = '  +0.44%  '
